# Trade #92 closed at 2026-02-17 15:54:28 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-ups for the MarketMaking
# strategy and appends the newly-closed trade (#92) to both the
# "All Trades" log and the per-strategy "MarketMaking" log.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - top level roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.78   # Current Capital
$summary.Range("B4").Value = -0.23     # Total P&L $
$summary.Range("B5").Value = -0.05     # Total P&L %
$summary.Range("B6").Value = 92        # Total Trades
$summary.Range("B7").Value = 32        # Winning Trades
$summary.Range("B9").Value = 34.78     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.78      # Capital
$status.Range("D4").Value = 92         # Trades
$status.Range("E4").Value = -0.23      # P&L $
$status.Range("F4").Value = -0.22      # P&L %
$status.Range("G4").Value = 34.78      # Win Rate %

# ---------------------------------------------------------------------
# "All Trades" and "MarketMaking" sheets - append closed trade #92
# ---------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($name in $tradeSheets) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A93").Value = 92
    # Pre-formatting as text ("@") before assigning keeps these as plain
    # text instead of being auto-parsed into date/time serial numbers,
    # matching the existing rows above them; ClearFormats() afterwards
    # drops the temporary text-format styling so the cell stays on the
    # sheet's default style, same as its neighbours.
    $ws.Range("B93").NumberFormat = "@"
    $ws.Range("B93").Value = "2026-02-17"
    $ws.Range("B93").ClearFormats()
    $ws.Range("C93").NumberFormat = "@"
    $ws.Range("C93").Value = "15:54:21"
    $ws.Range("C93").ClearFormats()
    $ws.Range("D93").Value = "MarketMaking"
    $ws.Range("E93").Value = "UP"
    $ws.Range("F93").Value = 0.72
    $ws.Range("G93").Value = 0.84
    $ws.Range("H93").Value = "CLOSED"
    $ws.Range("I93").Value = 16.6667
    $ws.Range("J93").Value = 0.12
    $ws.Range("K93").Value = 99.78
    $ws.Range("L93").Value = 0
    $ws.Range("M93").Value = 0
    $ws.Range("N93").Value = 0.6
    $ws.Range("O93").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P93").Value = "early_exit"
    $ws.Range("Q93").Value = 0.14
}
